$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shrink the second logo picture (top-right) ---
# Only the "to" (bottom-right) anchor's column offset shrinks; the "from" anchor is untouched,
# which shows up as a narrower picture. Width in points derived to hit the exact target colOff.
$logo = $ws.Shapes.Item(2)
$logo.Width = 170.1863

# --- Widen the data columns I:U ---
$ws.Columns.Item(9).ColumnWidth  = 16.666666666666668  # I
$ws.Columns.Item(10).ColumnWidth = 16.833333333333332  # J
$ws.Columns.Item(11).ColumnWidth = 17.333333333333332  # K
$ws.Columns.Item(12).ColumnWidth = 18.333333333333332  # L
$ws.Columns.Item(13).ColumnWidth = 19.166666666666668  # M
$ws.Columns.Item(14).ColumnWidth = 17.333333333333332  # N
$ws.Columns.Item(15).ColumnWidth = 18.166666666666668  # O
$ws.Columns.Item(16).ColumnWidth = 17.666666666666668  # P
$ws.Columns.Item(17).ColumnWidth = 18.5                # Q
$ws.Columns.Item(18).ColumnWidth = 17.333333333333332  # R
$ws.Columns.Item(19).ColumnWidth = 17.0                # S
$ws.Columns.Item(20).ColumnWidth = 19.5                # T
$ws.Columns.Item(21).ColumnWidth = 16.666666666666668  # U

# --- Sheet view: zoom out further in Page Break Preview and move the selection ---
$excel.ActiveWindow.Zoom = 50
$ws.Range("R18").Select()

# --- Print scale ---
$ws.PageSetup.Zoom = 38
